$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 34483960
$ws.Range("I137").Value = 43479056
$ws.Range("K137").Value = 130437168
$ws.Range("M137").Value = -130434618

$ws.Range("H138").Value = 4164745.8
$ws.Range("I138").Value = 1230703.2
$ws.Range("J138").Value = 6175108.5
$ws.Range("K138").Value = 3692109.6
$ws.Range("L138").Value = 18525325.5
$ws.Range("M138").Value = -3686969.6
$ws.Range("N138").Value = -18535605.5

$ws.Range("H141").Value = 2576.5
$ws.Range("I141").Value = 1702.326
$ws.Range("J141").Value = 5089.75
$ws.Range("K141").Value = 5106.978
$ws.Range("L141").Value = 15269.25
$ws.Range("M141").Value = 73.02199999999993
$ws.Range("N141").Value = -25629.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13345.574
$ws.Range("I32").Value = 2496.9634
$ws.Range("K32").Value = 2496.9634
$ws.Range("M32").Value = -2209.9634

$ws.Range("H61").Value = 2896.9048
$ws.Range("I61").Value = 1882.2258
$ws.Range("J61").Value = 5756.4546
$ws.Range("K61").Value = 1882.2258
$ws.Range("L61").Value = 5756.4546
$ws.Range("M61").Value = -1670.2258
$ws.Range("N61").Value = -6180.4546

$ws.Range("H74").Value = 8131.15
$ws.Range("I74").Value = 1926.9166
$ws.Range("K74").Value = 1926.9166
$ws.Range("M74").Value = -1052.9166

$ws.Range("H77").Value = 8131.15
$ws.Range("I77").Value = 1926.9166
$ws.Range("K77").Value = 9634.583000000001
$ws.Range("M77").Value = -5266.583000000001

$ws.Range("H132").Value = 3070.1
$ws.Range("I132").Value = 2731.0322
$ws.Range("K132").Value = 8193.096600000001
$ws.Range("M132").Value = -5663.096600000001

$ws.Range("H136").Value = 2896.9048
$ws.Range("I136").Value = 1882.2258
$ws.Range("J136").Value = 5756.4546
$ws.Range("K136").Value = 5646.6774
$ws.Range("L136").Value = 17269.3638
$ws.Range("M136").Value = -3096.6774
$ws.Range("N136").Value = -22369.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 8500
$ws.Range("J61").Value = 8500
$ws.Range("L61").Value = 8500
$ws.Range("N61").Value = -9126

$ws.Range("H107").Value = 1586.6364
$ws.Range("I107").Value = 1270.25
$ws.Range("J107").Value = 1767.4286
$ws.Range("K107").Value = 1270.25
$ws.Range("L107").Value = 1767.4286
$ws.Range("M107").Value = 649.75
$ws.Range("N107").Value = -5607.4286

$ws.Range("H134").Value = 17545892
$ws.Range("I134").Value = 21278166
$ws.Range("J134").Value = 4203.9
$ws.Range("K134").Value = 63834498
$ws.Range("L134").Value = 12611.7
$ws.Range("M134").Value = -63831963
$ws.Range("N134").Value = -17681.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 102142.2

$ws.Range("H31").Value = 1658.2742
$ws.Range("I31").Value = 1048.6757
$ws.Range("J31").Value = 2560.48
$ws.Range("K31").Value = 1048.6757
$ws.Range("L31").Value = 2560.48
$ws.Range("M31").Value = -753.6757
$ws.Range("N31").Value = -3150.48

$ws.Range("H34").Value = 1658.2742
$ws.Range("I34").Value = 1048.6757
$ws.Range("J34").Value = 2560.48
$ws.Range("K34").Value = 1048.6757
$ws.Range("L34").Value = 2560.48
$ws.Range("M34").Value = -846.6757
$ws.Range("N34").Value = -2964.48

$ws.Range("H58").Value = 1866.0264
$ws.Range("I58").Value = 1200.12
$ws.Range("K58").Value = 1200.12
$ws.Range("M58").Value = -997.1199999999999

$ws.Range("H113").Value = 102142.2

$ws.Range("H119").Value = 44380.5
$ws.Range("J119").Value = 44380.5
$ws.Range("L119").Value = 44380.5
$ws.Range("N119").Value = -54056.5

$ws.Range("H132").Value = 2123.5095
$ws.Range("I132").Value = 1684.2667
$ws.Range("K132").Value = 5052.800099999999
$ws.Range("M132").Value = -2522.800099999999

$ws.Range("H134").Value = 1944.5745
$ws.Range("I134").Value = 1411.2778
$ws.Range("J134").Value = 3689.9092
$ws.Range("K134").Value = 4233.8334
$ws.Range("L134").Value = 11069.7276
$ws.Range("M134").Value = -1698.8334
$ws.Range("N134").Value = -16139.7276

$ws.Range("H136").Value = 1866.0264
$ws.Range("I136").Value = 1200.12
$ws.Range("K136").Value = 3600.36
$ws.Range("M136").Value = -1050.36

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 231.8
$ws.Range("I13").Value = 86.666664
$ws.Range("J13").Value = 449.5
$ws.Range("K13").Value = 259.999992
$ws.Range("L13").Value = 1348.5
$ws.Range("M13").Value = -91.99999200000002
$ws.Range("N13").Value = -1684.5

$ws.Range("H17").Value = 220
$ws.Range("J17").Value = 220
$ws.Range("L17").Value = 660
$ws.Range("N17").Value = -998

$ws.Range("H92").Value = 806.0769
$ws.Range("I92").Value = 795
$ws.Range("J92").Value = 811
$ws.Range("K92").Value = 2385
$ws.Range("L92").Value = 2433
$ws.Range("M92").Value = -1137
$ws.Range("N92").Value = -4929

$ws.Range("H123").Value = 1614.4445
$ws.Range("I123").Value = 343.33334
$ws.Range("J123").Value = 2250
$ws.Range("K123").Value = 1030.00002
$ws.Range("L123").Value = 6750
$ws.Range("M123").Value = 1419.99998
$ws.Range("N123").Value = -11650

$ws.Range("H129").Value = 1633.8695
$ws.Range("I129").Value = 2033.3334
$ws.Range("J129").Value = 1492.8823
$ws.Range("K129").Value = 6100.0002
$ws.Range("L129").Value = 4478.6469
$ws.Range("M129").Value = -1100.0002
$ws.Range("N129").Value = -14478.6469

$ws.Range("H131").Value = 5849588
$ws.Range("J131").Value = 6804545
$ws.Range("L131").Value = 20413635
$ws.Range("N131").Value = -20423715

$ws.Range("H133").Value = 28934
$ws.Range("I133").Value = 2320
$ws.Range("J133").Value = 50709.09
$ws.Range("K133").Value = 6960
$ws.Range("L133").Value = 152127.27
$ws.Range("M133").Value = -1900
$ws.Range("N133").Value = -162247.27

$ws.Range("H134").Value = 3981.1155
$ws.Range("I134").Value = 2891.318
$ws.Range("J134").Value = 9975
$ws.Range("K134").Value = 8673.954000000002
$ws.Range("L134").Value = 29925
$ws.Range("M134").Value = -3603.954000000002
$ws.Range("N134").Value = -40065

$ws.Range("H136").Value = 2311.6
$ws.Range("I136").Value = 1310.8334
$ws.Range("J136").Value = 2833.739
$ws.Range("K136").Value = 3932.5002
$ws.Range("L136").Value = 8501.217000000001
$ws.Range("M136").Value = 1167.4998
$ws.Range("N136").Value = -18701.217

$ws.Range("H137").Value = 4043251.5
$ws.Range("I137").Value = 4350619.5
$ws.Range("J137").Value = 508516.5
$ws.Range("K137").Value = 13051858.5
$ws.Range("L137").Value = 1525549.5
$ws.Range("M137").Value = -13046758.5
$ws.Range("N137").Value = -1535749.5

$ws.Range("H138").Value = 1237.5
$ws.Range("I138").Value = 985.7143
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 2957.1429
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = 2182.8571
$ws.Range("N138").Value = -19280

$ws.Range("H139").Value = 3334.5557
$ws.Range("I139").Value = 3108.0908
$ws.Range("J139").Value = 3690.4285
$ws.Range("K139").Value = 9324.2724
$ws.Range("L139").Value = 11071.2855
$ws.Range("M139").Value = -4184.2724
$ws.Range("N139").Value = -21351.2855

$ws.Range("H140").Value = 6579.769
$ws.Range("I140").Value = 8111.107
$ws.Range("J140").Value = 2681.818
$ws.Range("K140").Value = 24333.321
$ws.Range("L140").Value = 8045.454000000001
$ws.Range("M140").Value = -19153.321
$ws.Range("N140").Value = -18405.454

$ws.Range("H141").Value = 3338.6667
$ws.Range("I141").Value = 3590
$ws.Range("J141").Value = 2333.3333
$ws.Range("K141").Value = 10770
$ws.Range("L141").Value = 6999.999899999999
$ws.Range("M141").Value = -5590
$ws.Range("N141").Value = -17359.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2751.9219
$ws.Range("I132").Value = 2573.8408
$ws.Range("J132").Value = 3143.7
$ws.Range("K132").Value = 7721.5224
$ws.Range("L132").Value = 9431.099999999999
$ws.Range("M132").Value = -5191.5224
$ws.Range("N132").Value = -14491.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 20776.25
$ws.Range("J106").Value = 20776.25
$ws.Range("L106").Value = 20776.25
$ws.Range("N106").Value = -23300.25

$ws.Range("H132").Value = 2700.7805
$ws.Range("I132").Value = 1574.3636
$ws.Range("J132").Value = 7347.25
$ws.Range("K132").Value = 4723.0908
$ws.Range("L132").Value = 22041.75
$ws.Range("M132").Value = -2193.0908
$ws.Range("N132").Value = -27101.75

$ws.Range("H136").Value = 3979.5
$ws.Range("I136").Value = 2341.1143
$ws.Range("J136").Value = 12171.429
$ws.Range("K136").Value = 7023.342900000001
$ws.Range("L136").Value = 36514.287
$ws.Range("M136").Value = -4473.342900000001
$ws.Range("N136").Value = -41614.287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2605.6
$ws.Range("I132").Value = 2780.5686
$ws.Range("J132").Value = 2135.9473
$ws.Range("K132").Value = 8341.7058
$ws.Range("L132").Value = 6407.841899999999
$ws.Range("M132").Value = -5811.7058
$ws.Range("N132").Value = -11467.8419

$ws.Range("H136").Value = 18715.508
$ws.Range("I136").Value = 26335.36
$ws.Range("J136").Value = 2205.8333
$ws.Range("K136").Value = 79006.08
$ws.Range("L136").Value = 6617.499899999999
$ws.Range("M136").Value = -76456.08
$ws.Range("N136").Value = -11717.4999
